$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (SamplesTab / B3): samples query text (content itself is unchanged by this edit,
# but must be re-applied so the workbook's shared-string table ends up in the expected state)
$ws.Range("B3").Value = 'Match (f)<--(g:genomic_info)
WHERE g.platform in [''Not specified in data'']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '''') as `Sample ID`,
 coalesce(p.participant_id,'''') as `Participant ID`,
 coalesce(s.study_name, '''') as `Study Name`,
 coalesce(s.phs_accession,'''') as `Accession`,
 coalesce(samp.sample_tumor_status,'''') as `Tumor`,
coalesce(samp.sample_type,'''') as `Analyte Type`
ORDER By samp.sample_id LIMIT 100'

# StatQuery column (C2:C4): new combined stats query text, same for every row
$statsQuery = 'CALL{
    MATCH (p:participant)-->(s:study)
  OPTIONAL  MATCH (samp:sample)-->(p)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE COALESCE(g.platform, "Not specified in data") in [''Not specified in data'']
    RETURN 
        count(distinct p) AS num_participants
}
WITH num_participants
CALL {
    MATCH (samp:sample)-->(p:participant)-->(s)
    MATCH (samp)<--(f:file)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE COALESCE(g.platforme, "Not specified in data") in [''Not specified in data'']
    RETURN 
        count(distinct samp) AS num_samples
}
WITH num_participants, num_samples
CALL {
    MATCH (f:file)-->(s:study)
    OPTIONAL MATCH (f)<--(g:genomic_info)
   OPTIONAL MATCH (samp:sample)<--(f)
    OPTIONAL MATCH (p:participant)<--(samp)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE COALESCE(g.platform, "Not specified in data") in [''Not specified in data'']
    RETURN 
        count(distinct s) AS num_studies,
        count(distinct f) AS num_files
}
RETURN 
    num_studies AS Studies,
    num_participants AS Participants,
    num_samples AS Samples,
    num_files AS `Files`'
$ws.Range("C2").Value = $statsQuery
$ws.Range("C3").Value = $statsQuery
$ws.Range("C4").Value = $statsQuery

# Row 2 (ParticipantsTab / B2): new participants query text
$ws.Range("B2").Value = 'MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE COALESCE(g.platform, "Not specified in data") in [''Not specified in data'']
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'''') as `Participant ID`,
coalesce(s.study_name, '''') as `Study Name`,
coalesce(s.phs_accession,'''') as `Accession`,
coalesce(p.gender,'''') as `Gender`,
coalesce(apoc.text.join(samp, '',''), '''') as `Samples`
ORDER BY p.participant_id limit 100'

# Row 4 (FilesTab / B4): new files query text
$ws.Range("B4").Value = 'MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE COALESCE(g.platform, "Not specified in data") in [''Not specified in data'']
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '''') as `File Name`,
    coalesce(s.study_name,'''') as `Study Name`,
    coalesce(s.phs_accession,'''') as `Accession`,
    coalesce(p.participant_id, '''') as `Participant ID`,
    coalesce(samp.sample_id, '''') as `Sample ID`,
    coalesce(f.file_type, '''') as `File Type`
ORDER BY f.file_name limit 100'

# Row heights grew to fit the much longer query text now stored in each row
$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(3).RowHeight = 409.5
$ws.Rows.Item(4).RowHeight = 409.5

# Update the active selection to B5
$ws.Range("B5").Select() | Out-Null
